$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying NATMI TPM recompute dropped every row whose "Target
# cluster" (column D) is "ECs" (rows 2, 5 and 8 of the original sheet),
# and refreshed every expression/specificity metric for the remaining
# Sending-cluster x Target-cluster combinations. Delete bottom-up so the
# row numbers of the rows still awaiting deletion don't shift.
$ws.Range("A8:T8").EntireRow.Delete()
$ws.Range("A5:T5").EntireRow.Delete()
$ws.Range("A2:T2").EntireRow.Delete()

# Remaining rows are now packed into rows 2-7 (dimension becomes A1:T7).
# Refresh them with the new TPM-derived values, column by column.
$data = @(
    @("ECs",   "Pspn", "Gfra4", "FAPs",  2, 0.6666666666666666, 0.6034996666666667, 1.810499,   0.4381949890118401, 0.4381949890118402, 2, 0.6666666666666666, 0.3037343333333333, 0.911203, 0.5768794419634096, 0.5768794419634096, 0.1833035689218889,  1.649732120297,     0.2527856807323127, 0.2527856807323127),
    @("ECs",   "Pspn", "Gfra4", "MuSCs", 2, 0.6666666666666666, 0.6034996666666667, 1.810499,   0.4381949890118401, 0.4381949890118402, 3, 1,                  0.2227783333333333, 0.668335, 0.4231205580365905, 0.4231205580365905, 0.1344466499072222,  1.210019849165,     0.1854093082795274, 0.1854093082795275),
    @("FAPs",  "Pspn", "Gfra4", "FAPs",  1, 0.3333333333333333, 0.4194516666666666, 1.258355,   0.3045596022963802, 0.3045596022963802, 2, 0.6666666666666666, 0.3037343333333333, 0.911203, 0.5768794419634096, 0.5768794419634096, 0.1274018723405555,  1.146616851065,     0.1756941734173338, 0.1756941734173338),
    @("FAPs",  "Pspn", "Gfra4", "MuSCs", 1, 0.3333333333333333, 0.4194516666666666, 1.258355,   0.3045596022963802, 0.3045596022963802, 3, 1,                  0.2227783333333333, 0.668335, 0.4231205580365905, 0.4231205580365905, 0.09344474321388888, 0.841002688925,     0.1288654288790465, 0.1288654288790465),
    @("MuSCs", "Pspn", "Gfra4", "FAPs",  2, 0.6666666666666666, 0.3542886666666667, 1.062866,   0.2572454086917797, 0.2572454086917798, 2, 0.6666666666666666, 0.3037343333333333, 0.911203, 0.5768794419634096, 0.5768794419634096, 0.1076096319775556,  0.9684866877980001, 0.1483995878137631, 0.1483995878137631),
    @("MuSCs", "Pspn", "Gfra4", "MuSCs", 2, 0.6666666666666666, 0.3542886666666667, 1.062866,   0.2572454086917797, 0.2572454086917798, 3, 1,                  0.2227783333333333, 0.668335, 0.4231205580365905, 0.4231205580365905, 0.07892783867888889, 0.7103505481100001, 0.1088458208780166, 0.1088458208780166)
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $columns.Length; $j++) {
        $ws.Range($columns[$j] + $rowNum).Value = $rowValues[$j]
    }
}
